$d = $word.ActiveDocument

# Move to end of document and add a new paragraph
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd
$end.InsertParagraphAfter()

Write-Output $d.Paragraphs.Count
Write-Output ("[" + $d.Content.Text + "]")

$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Para without auto spacing"

Write-Output ("[" + $d.Content.Text + "]")

$p2.SpaceBefore = 20
$p2.SpaceAfter = 20
$p2.SpaceBeforeAuto = $false
$p2.SpaceAfterAuto = $false

Write-Output $p2.SpaceBefore
Write-Output $p2.SpaceAfter
Write-Output $p2.SpaceBeforeAuto
Write-Output $p2.SpaceAfterAuto

